# Update the ClientAdminCreation sheet: rename the automation supplier
# alias used by the end-user integrated test rows, flag the skipped /
# failing run-mode markers, and leave the sheet active with A4 selected
# (matching the "end user integrated files" commit).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ClientAdminCreation")

# clientSupplier / clientSupplierAlias values used by rows 3-5
$ws.Range("J3:J5").Value = "Auto Supplier"
$ws.Range("K3:K5").Value = "auto"

# Results column: header row stays blank, data rows now carry a status
$ws.Range("N2").Value = "SKIP"
$ws.Range("N3:N5").Value = "FAIL"

# Make this the active sheet/selection, as it is the one just edited
$ws.Activate()
$ws.Range("A4").Select()
